$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> (new DAMSLTag, new DialogAct)
# Derived from the canonical OOXML diff (97 row re-annotations after SGNN re-run).
$updates = @(
    @{ Row = 25; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 26; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 30; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 34; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 46; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 53; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 56; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 62; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 72; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 74; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 76; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 82; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 85; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 96; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 107; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 110; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 130; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 133; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 143; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 160; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 161; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 165; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 169; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 173; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 180; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 189; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 207; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 214; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 220; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 230; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 236; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 238; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 252; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 275; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 276; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 286; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 290; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 310; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 331; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 344; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 349; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 355; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 357; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 362; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 363; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 369; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 371; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 380; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 389; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 393; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 400; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 410; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 415; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 426; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 429; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 430; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 432; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 434; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 437; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 454; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 460; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 462; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 465; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 467; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 470; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 476; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 491; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 503; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 505; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 526; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 540; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 544; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 554; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 565; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 567; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 575; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 578; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 588; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 592; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 594; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 598; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 599; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 604; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 605; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 614; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 633; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 652; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 660; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 661; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 679; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 682; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 688; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 693; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 702; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 706; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 715; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 730; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

Write-Output ("Updated " + $updates.Count + " rows.")
